# Template for Shiny app
# - drop the "Date" row from the metadata sheet (rows below shift up)
# - tidy up the saved selections on metadata / key
# - make the Animal sheet the active tab
# - turn the Animal sheet's AnimalID column into a formula ("Animal_<n>")

$wb = $excel.ActiveWorkbook

# --- metadata sheet: remove the "Date" row (old row 3); rows 4-6 shift up to 3-5 ---
$wsMeta = $wb.Worksheets.Item("metadata")
$wsMeta.Rows.Item(3).Delete() | Out-Null
$wsMeta.Range("B2").Select() | Out-Null

# --- key sheet: just a saved-selection change (loses tabSelected once Animal
#     is activated below) ---
$wsKey = $wb.Worksheets.Item("key")
$wsKey.Range("B5").Select() | Out-Null

# --- Animal sheet: column A becomes ="Animal_"&ROW()-1 (A2 standalone,
#     A3:A11 as a shared-formula block), and the sheet becomes the active tab ---
$wsAnimal = $wb.Worksheets.Item("Animal")
$wsAnimal.Range("A2").Formula = "=""Animal_""&ROW()-1"
$wsAnimal.Range("A3:A11").Formula = "=""Animal_""&ROW()-1"

$wsAnimal.Activate() | Out-Null
$wsAnimal.Range("A4").Select() | Out-Null
